$d = $word.ActiveDocument

$replacements = @(
    @('2025-08-09 Saturday', '2025-08-10 Sunday'),
    @('870×5=4350', '387×9=3483'),
    @('563×9=5067', '259×3=777'),
    @('945×2=1890', '682×5=3410'),
    @('378×4=1512', '295×5=1475'),
    @('304×7=2128', '176×3=528'),
    @('654×5=3270', '681×9=6129'),
    @('824×9=7416', '196×8=1568'),
    @('578×5=2890', '972×6=5832'),
    @('196×7=1372', '116×3=348'),
    @('438×4=1752', '399×7=2793'),
    @('921×6=5526', '313×4=1252'),
    @('494×8=3952', '979×2=1958'),
    @('855×7=5985', '417×9=3753'),
    @('729×2=1458', '401×7=2807'),
    @('138×5=690', '281×6=1686'),
    @('816×3=2448', '982×5=4910'),
    @('948×8=7584', '853×3=2559'),
    @('261×3=783', '227×4=908'),
    @('423×2=846', '948×2=1896'),
    @('343×6=2058', '599×6=3594'),
    @('816×2=1632', '837×8=6696'),
    @('750×6=4500', '366×3=1098'),
    @('599×7=4193', '935×6=5610'),
    @('660×3=1980', '855×2=1710'),
    @('621×9=5589', '500×5=2500'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done"
